# 11 May - Noche
# Correct the Materia (E) / Docente (F) pairing on the "Blancos" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blancos")

$rows = @{
    3  = @("TEMAS DE FILOSOFÍA", "Hernández Mendoza Delfina")
    4  = @("PROBABILIDAD Y ESTADÍSTICA", "Velasco Sanchez David")
    5  = @("REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS", "Bautista Sarao Eutiquio")
    6  = @("TEMAS DE FÍSICA", "Duran Amezcua Maria Angelica")
    7  = @("REALIZA MANTENIMIENTO EN EL SISTEMA DE DISTRIBUCIÓN DE ENERGÍA ELÉCTRICA", "Bautista Sarao Eutiquio")
    16 = @("REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS", "Bautista Sarao Eutiquio")
    19 = @("TEMAS DE FÍSICA", "Duran Amezcua Maria Angelica")
    24 = @("TEMAS DE FILOSOFÍA", "Hernández Mendoza Delfina")
    25 = @("PROBABILIDAD Y ESTADÍSTICA", "Velasco Sanchez David")
    26 = @("TEMAS DE FÍSICA", "Duran Amezcua Maria Angelica")
    27 = @("REALIZA MANTENIMIENTO EN EL SISTEMA DE DISTRIBUCIÓN DE ENERGÍA ELÉCTRICA", "Bautista Sarao Eutiquio")
    29 = @("PROBABILIDAD Y ESTADÍSTICA", "Velasco Sanchez David")
    30 = @("TEMAS DE FILOSOFÍA", "Hernández Mendoza Delfina")
    33 = @("REALIZA MANTENIMIENTO EN EL SISTEMA DE DISTRIBUCIÓN DE ENERGÍA ELÉCTRICA", "Bautista Sarao Eutiquio")
    34 = @("REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS", "Bautista Sarao Eutiquio")
    45 = @("TEMAS DE FÍSICA", "Duran Amezcua Maria Angelica")
    47 = @("REALIZA MANTENIMIENTO EN EL SISTEMA DE DISTRIBUCIÓN DE ENERGÍA ELÉCTRICA", "Bautista Sarao Eutiquio")
    48 = @("TEMAS DE FILOSOFÍA", "Hernández Mendoza Delfina")
    49 = @("PROBABILIDAD Y ESTADÍSTICA", "Velasco Sanchez David")
    50 = @("TEMAS DE FÍSICA", "Duran Amezcua Maria Angelica")
    51 = @("REALIZA MANTENIMIENTO EN EL SISTEMA DE DISTRIBUCIÓN DE ENERGÍA ELÉCTRICA", "Bautista Sarao Eutiquio")
    53 = @("TEMAS DE FILOSOFÍA", "Hernández Mendoza Delfina")
    54 = @("PROBABILIDAD Y ESTADÍSTICA", "Velasco Sanchez David")
}

foreach ($r in $rows.Keys) {
    $pair = $rows[$r]
    $ws.Cells.Item($r, 5).Value = $pair[0]
    $ws.Cells.Item($r, 6).Value = $pair[1]
}
